# Regenerate orders with updated distance/size codes.
#
# The stimulus-naming scheme encodes a viewing Distance (D51/D64/D80) and a
# face Size (S20/S25/S30) inside several columns (Condition, Filename_Left,
# Filename_Right, Distance, Size). This run updates the calibration so that:
#   D80 -> D86
#   D64 -> D69
#   D51 -> D55
#   S30 -> S31
# everywhere those codes appear (e.g. "Face11_D80_S20" -> "Face11_D86_S20",
# "Fixation_D64_l.png" -> "Fixation_D69_l.png", "D51" -> "D55", etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

$updated = 0
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val.Replace("D80", "D86").Replace("D64", "D69").Replace("D51", "D55").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
                $updated = $updated + 1
            }
        }
    }
}

"Updated $updated cell(s) with new distance/size codes."
